$d = $word.ActiveDocument

# --- Update the letter's date line -----------------------------------
# Before: "31 March 2017"
# After : "03 April 2017"
#
# Locate the existing date text with Find rather than a hard-coded
# paragraph index, then capture where it starts so we can re-derive
# positions inside the replaced text afterwards.
$dateRange = $d.Content
$found = $dateRange.Find.Execute("31 March 2017", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the original date text '31 March 2017'"
}

$dateStart = $dateRange.Start

# Replace the whole date string in one go.
$dateRange.Text = "03 April 2017"

# --- Relocate the "_GoBack" bookmark ----------------------------------
# In the original file the (hidden) "_GoBack" bookmark sat further down,
# right after the "April 07" run inside the main paragraph (between
# "April 07" and ", 2017"). The edit moves that same bookmark so it now
# sits right after "03" in the new date line. Word only ever keeps one
# bookmark per name, so re-adding "_GoBack" at the new location removes
# it from its old spot automatically.
$bmPos = $dateStart + 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
